# CSV file upload was integrated - update student address values to remove commas
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Address column (F) for the two student rows: strip commas from the address
$ws.Range("F2").Value = "Sango Ota Ogun State Nigeria"
$ws.Range("F3").Value = "Sango Ota Ogun State Nigeria"

# Reflect the final selection/view state left by the edit (selection ends on F3,
# and the view scrolls back so column A is the left-most visible column)
$ws.Range("F3").Select()
$excel.ActiveWindow.ScrollColumn = 1
